{"js": "// Edit: \"Les librairies suivantes sont \u00e0 installer sur l'ordinateur qui va\n// analyser l'image :\" ->\n// \"Les librairies suivantes sont \u00e0 installer sur l'ordinateur qui va\n// analyser l'image et sur l'AlphaBot 2 :\"\n//\n// (the trailing non-breaking-space + colon that was already in the\n// document is left untouched; we only insert \" et sur l'AlphaBot 2\"\n// right after \"image\").\n\nconst body = context.document.body;\n\n// This exact sentence only occurs once in the document (the near-identical\n// sentence in the \"Serveur / AlphaBot2\" section ends in \"l'AlphaBot 2\"\n// instead of \"l'ordinateur qui va analyser l'image\", so this search string\n// is unique).\nconst target =\n  \"Les librairies suivantes sont \\u00e0 installer sur l\\u2019ordinateur qui va analyser l\\u2019image\";\n\nconst results = body.search(target, { matchCase: false, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found\");\n}\n\nconst found = results.items[0];\n\n// Insert the new clause right after \"...l'image\", before the existing\n// (non-breaking-space + colon) that closes the sentence.\nfound.insertText(\" et sur l\\u2019AlphaBot 2\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Edit: \"Les librairies suivantes sont \u00e0 installer sur l'ordinateur qui va\n# analyser l'image :\" ->\n# \"Les librairies suivantes sont \u00e0 installer sur l'ordinateur qui va\n# analyser l'image et sur l'AlphaBot 2 :\"\n#\n# (the trailing non-breaking-space + colon that was already in the\n# document is left untouched; we only insert \" et sur l'AlphaBot 2\"\n# right after \"image\".)\n\n$d = $word.ActiveDocument\n\n# This exact sentence only occurs once in the document (the near-identical\n# sentence in the \"Serveur / AlphaBot2\" section ends in \"l'AlphaBot 2\"\n# instead of \"l'ordinateur qui va analyser l'image\", so this search string\n# is unique).\n$range = $d.Content\n$range.Find.MatchCase = $false\n$found = $range.Find.Execute(\"Les librairies suivantes sont \u00e0 installer sur l\u2019ordinateur qui va analyser l\u2019image\")\n\nif (-not $range.Find.Found) {\n    throw \"Target sentence not found\"\n}\n\n# Collapse to the end of the matched range, then insert the new clause\n# right before the existing (non-breaking-space + colon) that closes the\n# sentence.\n$range.Collapse($wdCollapseEnd)\n$range.InsertAfter(\" et sur l\u2019AlphaBot 2\")\n"}
